$d = $word.ActiveDocument

# Locate the run containing the buggy formula text.
$rng = $d.Content
$found = $rng.Find.Execute("MWTPa = 0.513*()")

if ($found) {
    $target = $d.Range($rng.Start, $rng.End)
    $target.Delete()

    # 1) Plain run: "MWTPa = 0.513"
    $target.InsertAfter("MWTPa = 0.513")
    $target.Collapse(0)

    # 2) Italic run: "(Q0-Q) + 1.662"
    $target.InsertAfter("(Q0-Q) + 1.662")
    $target.Italic = 1
    $target.Collapse(0)

    # 3) Italic run: " " (space)
    $target.InsertAfter(" ")
    $target.Italic = 1
    $target.Collapse(0)

    # 4) Italic run: "MWTPa = 0.513"
    $target.InsertAfter("MWTPa = 0.513")
    $target.Italic = 1
    $target.Collapse(0)

    # 5) Plain run: "(170 - Q) + 1.662"
    $target.InsertAfter("(170 - Q) + 1.662")
}
